$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "NetNratio"
$ws.Range("N2").Formula = "=D2/I2"
$ws.Range("N3:N19").Formula = "=D3/I3"
$ws.Columns("N").AutoFit() | Out-Null

$ws.Range("N2").Select()
